$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: 2nd place
$ws.Range("C3").Value = "felipe"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2446"
$ws.Range("D3").Style = "Normal"

# Row 4: 3rd place
$ws.Range("C4").Value = "matheus"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2430"
$ws.Range("D4").Style = "Normal"

# Row 5: 4th place
$ws.Range("C5").Value = "Diegowl"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2344"
$ws.Range("D5").Style = "Normal"

# Row 6: 5th place
$ws.Range("C6").Value = "Esther linda"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2310"
$ws.Range("D6").Style = "Normal"

# Row 7: 6th place
$ws.Range("C7").Value = "bona"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1951"
$ws.Range("D7").Style = "Normal"

# Row 8: 7th place
$ws.Range("C8").Value = "oi"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1944"
$ws.Range("D8").Style = "Normal"

# Row 9: 8th place
$ws.Range("C9").Value = "last dance"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0"
$ws.Range("D9").Style = "Normal"

# Remove rows 10-12 (previously 9th, 10th, 11th place entries)
$ws.Range("A10:D12").Delete()
